$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1971.4762
$ws.Range("J40").Value = 2026.3158
$ws.Range("L40").Value = 2026.3158
$ws.Range("N40").Value = -2376.3158
$ws.Range("H64").Value = 4550
$ws.Range("I64").Value = 6600
$ws.Range("J64").Value = 3525
$ws.Range("K64").Value = 6600
$ws.Range("L64").Value = 3525
$ws.Range("M64").Value = -6352
$ws.Range("N64").Value = -4021
$ws.Range("H67").Value = 4550
$ws.Range("I67").Value = 6600
$ws.Range("J67").Value = 3525
$ws.Range("K67").Value = 6600
$ws.Range("L67").Value = 3525
$ws.Range("M67").Value = -5742
$ws.Range("N67").Value = -5241
$ws.Range("H131").Value = 2272.1365
$ws.Range("I131").Value = 943.8182
$ws.Range("J131").Value = 3600.4546
$ws.Range("K131").Value = 2831.4546
$ws.Range("L131").Value = 10801.3638
$ws.Range("M131").Value = 2208.5454
$ws.Range("N131").Value = -20881.3638
$ws.Range("H137").Value = 1416
$ws.Range("I137").Value = 1158.1875
$ws.Range("J137").Value = 2241
$ws.Range("K137").Value = 3474.5625
$ws.Range("L137").Value = 6723
$ws.Range("M137").Value = -924.5625
$ws.Range("N137").Value = -11823
$ws.Range("H138").Value = 4152.771
$ws.Range("I138").Value = 2529.2222
$ws.Range("J138").Value = 4527.436
$ws.Range("K138").Value = 7587.6666
$ws.Range("L138").Value = 13582.308
$ws.Range("M138").Value = -2447.6666
$ws.Range("N138").Value = -23862.308
$ws.Range("H141").Value = 766.6667
$ws.Range("I141").Value = 650
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 1950
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 3230
$ws.Range("N141").Value = -13360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7618.615
$ws.Range("I86").Value = 8955.889
$ws.Range("J86").Value = 4609.75
$ws.Range("K86").Value = 8955.889
$ws.Range("L86").Value = 4609.75
$ws.Range("M86").Value = -7832.888999999999
$ws.Range("N86").Value = -6855.75
$ws.Range("H89").Value = 7618.615
$ws.Range("I89").Value = 8955.889
$ws.Range("J89").Value = 4609.75
$ws.Range("K89").Value = 44779.44499999999
$ws.Range("L89").Value = 23048.75
$ws.Range("M89").Value = -39163.44499999999
$ws.Range("N89").Value = -34280.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3414.907
$ws.Range("I31").Value = 2788.75
$ws.Range("J31").Value = 4205.8423
$ws.Range("K31").Value = 2788.75
$ws.Range("L31").Value = 4205.8423
$ws.Range("M31").Value = -2493.75
$ws.Range("N31").Value = -4795.8423
$ws.Range("H34").Value = 3414.907
$ws.Range("I34").Value = 2788.75
$ws.Range("J34").Value = 4205.8423
$ws.Range("K34").Value = 2788.75
$ws.Range("L34").Value = 4205.8423
$ws.Range("M34").Value = -2586.75
$ws.Range("N34").Value = -4609.8423
$ws.Range("H62").Value = 6538.5
$ws.Range("I62").Value = 6626.4287
$ws.Range("J62").Value = 6333.3335
$ws.Range("K62").Value = 6626.4287
$ws.Range("L62").Value = 6333.3335
$ws.Range("M62").Value = -6002.4287
$ws.Range("N62").Value = -7581.3335
$ws.Range("H65").Value = 6538.5
$ws.Range("I65").Value = 6626.4287
$ws.Range("J65").Value = 6333.3335
$ws.Range("K65").Value = 33132.14350000001
$ws.Range("L65").Value = 31666.6675
$ws.Range("M65").Value = -30012.14350000001
$ws.Range("N65").Value = -37906.6675

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 13907.143
$ws.Range("J82").Value = 13907.143
$ws.Range("L82").Value = 41721.429
$ws.Range("N82").Value = -42533.429
$ws.Range("H85").Value = 13907.143
$ws.Range("J85").Value = 13907.143
$ws.Range("L85").Value = 41721.429
$ws.Range("N85").Value = -44529.429
$ws.Range("H113").Value = 638.2727
$ws.Range("I113").Value = 729.55554
$ws.Range("J113").Value = 604.0417
$ws.Range("K113").Value = 2188.66662
$ws.Range("L113").Value = 1812.1251
$ws.Range("M113").Value = -18.66661999999997
$ws.Range("N113").Value = -6152.1251
$ws.Range("H140").Value = 132803.88
$ws.Range("I140").Value = 152374.45
$ws.Range("J140").Value = 2333.3333
$ws.Range("K140").Value = 457123.35
$ws.Range("L140").Value = 6999.999899999999
$ws.Range("M140").Value = -451943.35
$ws.Range("N140").Value = -17359.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 46.77778
$ws.Range("I2").Value = 37
$ws.Range("J2").Value = 81
$ws.Range("K2").Value = 37
$ws.Range("L2").Value = 81
$ws.Range("M2").Value = 76
$ws.Range("N2").Value = -307
$ws.Range("H11").Value = 11821127
$ws.Range("J11").Value = 8100
$ws.Range("L11").Value = 8100
$ws.Range("N11").Value = -8378
$ws.Range("H70").Value = 13041.546
$ws.Range("I70").Value = 17619.965
$ws.Range("J70").Value = 4189.933
$ws.Range("K70").Value = 17619.965
$ws.Range("L70").Value = 4189.933
$ws.Range("M70").Value = -17349.965
$ws.Range("N70").Value = -4729.933
$ws.Range("H73").Value = 13041.546
$ws.Range("I73").Value = 17619.965
$ws.Range("J73").Value = 4189.933
$ws.Range("K73").Value = 17619.965
$ws.Range("L73").Value = 4189.933
$ws.Range("M73").Value = -16683.965
$ws.Range("N73").Value = -6061.933
$ws.Range("H132").Value = 3076.353
$ws.Range("J132").Value = 4409.778
$ws.Range("L132").Value = 13229.334
$ws.Range("N132").Value = -18289.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2573.2917
$ws.Range("I68").Value = 2822.4375
$ws.Range("J68").Value = 2075
$ws.Range("K68").Value = 2822.4375
$ws.Range("L68").Value = 2075
$ws.Range("M68").Value = -2073.4375
$ws.Range("N68").Value = -3573
$ws.Range("H71").Value = 2573.2917
$ws.Range("I71").Value = 2822.4375
$ws.Range("J71").Value = 2075
$ws.Range("K71").Value = 14112.1875
$ws.Range("L71").Value = 10375
$ws.Range("M71").Value = -10368.1875
$ws.Range("N71").Value = -17863
$ws.Range("H82").Value = 2918.1875
$ws.Range("I82").Value = 2515.8333
$ws.Range("J82").Value = 4125.25
$ws.Range("K82").Value = 2515.8333
$ws.Range("L82").Value = 4125.25
$ws.Range("M82").Value = -2154.8333
$ws.Range("N82").Value = -4847.25
$ws.Range("H85").Value = 2918.1875
$ws.Range("I85").Value = 2515.8333
$ws.Range("J85").Value = 4125.25
$ws.Range("K85").Value = 2515.8333
$ws.Range("L85").Value = 4125.25
$ws.Range("M85").Value = -1267.8333
$ws.Range("N85").Value = -6621.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 9666.667
$ws.Range("I54").Value = 9666.667
$ws.Range("K54").Value = 9666.667
$ws.Range("M54").Value = -9146.667
$ws.Range("H62").Value = 4783.3335
$ws.Range("I62").Value = 5075
$ws.Range("K62").Value = 5075
$ws.Range("M62").Value = -4451
$ws.Range("H65").Value = 4783.3335
$ws.Range("I65").Value = 5075
$ws.Range("K65").Value = 25375
$ws.Range("M65").Value = -22255
